# Update D13 value from 1.453125 to 1.925 and its number format from
# the built-in "mmm-yy" (numFmtId 17) to "0.00" (numFmtId 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = 1.925
$ws.Range("D13").NumberFormat = "0.00"

# Move the active selection to D13 (was E19).
$ws.Range("D13").Select()
